# Updates cryptos list values (Price column D, Volume(1h) column E)
# to match the latest scrape, as produced by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.448.55'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '1.571.27'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '290.32'
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3694'
$ws.Range('E7').Value = '  -1.67%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '49.92'
$ws.Range('E8').Value = '  +1.48%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3383'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  +2.36%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07556'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.21'
$ws.Range('E13').Value = '  +2.02%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.031'
$ws.Range('E14').Value = '  +2.24%  '
$ws.Range('E15').Value = '  +1.89%  '
$ws.Range('D16').Value = '1.571.04'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '90.39'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06775'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.359'
$ws.Range('E21').Value = '  +3.16%  '
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('E23').Value = '  +3.03%  '
$ws.Range('D24').Value = '22.440.60'
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.358'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.670'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.00'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '149.25'
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.050'
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.92'
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('D31').Value = '1.747.92'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.065'
$ws.Range('E32').Value = '  +8.53%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.228'
$ws.Range('E33').Value = '  +4.71%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.016'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.845'
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08370'
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('E37').Value = '  +1.23%  '
$ws.Range('E38').Value = '  -3.42%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.2305'
$ws.Range('E39').Value = '  +2.07%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.06573'
$ws.Range('E40').Value = '  +2.73%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.428'
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.33'
$ws.Range('E42').Value = '  +3.30%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6262'
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.802'
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5888'
$ws.Range('E47').Value = '  +1.56%  '
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '128.13'
$ws.Range('E49').Value = '  +3.47%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.243'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07306'
